$d = $word.ActiveDocument

# Build a minimal single-part OOXML "package" wrapper so the new
# paragraph's <w:t> keeps an explicit xml:space="preserve" attribute
# (matching the rest of the document) instead of the engine's default
# whitespace-driven serialization.
function New-OoxmlDocPart([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + $bodyXml + '</pkg:xmlData>' + `
        '</pkg:part></pkg:package>'
}

# Insert a new "Author"-styled paragraph immediately after $para,
# containing the given name, and return the new paragraph.
function Add-AuthorAfter($para, [string]$name) {
    $null = $para.Range.InsertParagraphAfter()
    $newPara = $para.Next()
    $bodyXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:pPr><w:pStyle w:val="Author"/></w:pPr>' + `
        '<w:r><w:t xml:space="preserve">' + $name + '</w:t></w:r></w:p>'
    $newPara.Range.InsertXML((New-OoxmlDocPart $bodyXml))
    return $para.Next()
}

# Locate the "Colin Carlson" author paragraph and add the two new
# authors ("Elodie Eiffener", then "Gabriel Munoz Acevedo") right
# after it, before "Andrea Paz Velez".
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Colin Carlson") {
        $after = Add-AuthorAfter $p "Elodie Eiffener"
        $after = Add-AuthorAfter $after "Gabriel Munoz Acevedo"
        break
    }
}
